# Regenerate save_data: update column G ("K" = strikeouts) values for
# Josh Osich's 2021 save data, replacing the old "Strike#" derived values
# with the recalculated K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 3
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 2
    20 = 3
    21 = 1
    22 = 3
    23 = 2
    25 = 2
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
